{"js": "// The edit moves the existing \"_GoBack\" bookmark (a zero-length marker Word\n// drops at the last edited location) from the end of the \"Compl\u00e9t\u00e9.\"\n// paragraph to right before the word \"correct\" in the preceding paragraph\n// (\"...si le mot de passe est [BOOKMARK]correct, sinon retourne...\").\n// The visible text/runs are otherwise unchanged.\n\n// Remove the bookmark from its current location (no-op if already absent).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Locate the word \"correct\" that follows \"si le mot de passe est \".\nconst results = context.document.body.search(\"correct\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"correct\" to re-anchor the _GoBack bookmark.');\n}\n\n// Re-insert the bookmark as a collapsed range right before \"correct\".\nconst target = results.items[0].getRange(\"Start\");\ntarget.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The edit moves the existing \"_GoBack\" bookmark (a zero-length marker Word\n# drops at the last edited location) from the end of the \"Compl\u00e9t\u00e9.\"\n# paragraph to right before the word \"correct\" in the preceding paragraph\n# (\"...si le mot de passe est [BOOKMARK]correct, sinon retourne...\").\n# The visible text/runs are otherwise unchanged.\n\n$d = $word.ActiveDocument\n\n# Remove the bookmark from its current location (no-op if already absent).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Locate the word \"correct\" that follows \"si le mot de passe est \".\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 1  # wdFindContinue\n$rng.Find.MatchCase = $true\n$found = $rng.Find.Execute(\"correct\")\n\nif (-not $found) {\n    throw 'Could not find \"correct\" to re-anchor the _GoBack bookmark.'\n}\n\n# Collapse the found range to its start and re-insert the bookmark there.\n$target = $d.Range($rng.Start, $rng.Start)\n$d.Bookmarks.Add(\"_GoBack\", $target) | Out-Null\n"}
